# Apply "multi-industry template" restoration edits:
# Replace "Product Development"/"Product"-themed copy with "AI/ML"/"AI"-themed
# copy across the three worksheets of the Change Management Plan workbook.

$wb = $excel.ActiveWorkbook

# --- Sheet: "Change Management Overview" ---
$ws1 = $wb.Worksheets.Item("Change Management Overview")

$ws1.Range("A2").Value = "PRODUCT Change Management Plan Project"
$ws1.Range("B6").Value = "Enterprise AI/ML Implementation"
$ws1.Range("A15").Value = "1. Achieve 95% user adoption of new AI/ML systems within 6 months of go-live"
$ws1.Range("A17").Value = "3. Build organizational capability and confidence in AI/ML technologies"
$ws1.Range("A20").Value = "6. Create positive stakeholder sentiment and enthusiasm for AI/ML transformation"

# --- Sheet: "Change Impact Assessment" ---
$ws2 = $wb.Worksheets.Item("Change Impact Assessment")

$ws2.Range("G4").Value = "AI/ML automation"
$ws2.Range("G5").Value = "AI-powered insights"
$ws2.Range("G7").Value = "New AI interface"
$ws2.Range("G11").Value = "AI-enhanced CRM"
$ws2.Range("G12").Value = "AI-assisted support"
$ws2.Range("G13").Value = "AI-powered testing"
